$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert the "end"/"start" timestamp cells in column C from a raw date
# serial to an explicit text timestamp (timezone-aware string).
# C57 ("end") is set first so its shared string is interned before C56's,
# matching the authoring order of the workbook being replicated.
$ws.Range("C57").ClearFormats()
$ws.Range("C57").Value = "2020-12-31 23:00:00-07:00"
$ws.Range("C56").ClearFormats()
$ws.Range("C56").Value = "2018-05-24 00:00:00-07:00"

# Shorten manufacturer names in header row 1
$ws.Range("B1").Value = "LG"
$ws.Range("D1").Value = "Can270"
$ws.Range("E1").Value = "Can275"
$ws.Range("H1").Value = "MissionSolar"
$ws.Range("I1").Value = "HanQPlus"
$ws.Range("J1").Value = "HanQPeak"

# Move the active selection to I2
$null = $ws.Range("I2").Select()
